# Add BiggestTestCase2 and BiggestTestCase3.
#
# Inserts a new list paragraph right after the "... BiggestTestCase1."
# paragraph, containing the instructions for running the tests via the
# Test Explorer. The new paragraph inherits the same ListParagraph /
# numbered-list formatting as its neighbours.

$d = $word.ActiveDocument

# Locate the paragraph that ends with "BiggestTestCase1." using Find,
# then insert a brand-new paragraph right after it.
$searchRange = $d.Content.Duplicate
$found = $searchRange.Find.Execute("BiggestTestCase1.", $true, $false, $false,
                                    $false, $false, $true, 1, $false, "", 0)

$anchorPara = $searchRange.Paragraphs(1)
$newParaRange = $anchorPara.Range.InsertParagraphAfter()

# The freshly inserted paragraph is now the paragraph right after the anchor.
$newPara = $anchorPara.Next()
$newPara.Range.Text = "Right-click, then select Run Tests. This will open the Test Explorer."
